$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-coerce cells whose new numeric-looking value would otherwise be
# auto-converted to a number by Excel; applied before assignment, then the
# cell style is reset back to Normal so no stray formatting is introduced.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D17", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value = "26.563.70"
$ws.Range("E2").Value = "  -7.52%  "
$ws.Range("D3").Value = "1.686.94"
$ws.Range("E3").Value = "  -6.56%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "217.37"
$ws.Range("E5").Value = "  -6.30%  "
$ws.Range("D6").Value = "0.5011"
$ws.Range("E6").Value = "  -15.76%  "
$ws.Range("D8").Value = "0.2621"
$ws.Range("E8").Value = "  -5.96%  "
$ws.Range("D9").Value = "21.94"
$ws.Range("E9").Value = "  -6.06%  "
$ws.Range("D10").Value = "0.06212"
$ws.Range("E10").Value = "  -9.27%  "
$ws.Range("D11").Value = "0.07288"
$ws.Range("E11").Value = "  -3.04%  "
$ws.Range("D12").Value = "1.637.06"
$ws.Range("E12").Value = "  -9.27%  "
$ws.Range("D13").Value = "4.451"
$ws.Range("E13").Value = "  -7.54%  "
$ws.Range("D14").Value = "0.5775"
$ws.Range("E14").Value = "  -7.50%  "
$ws.Range("D15").Value = "1.916.10"
$ws.Range("E15").Value = "  -6.58%  "
$ws.Range("D16").Value = "0.000008199"
$ws.Range("E16").Value = "  -12.19%  "
$ws.Range("D17").Value = "64.82"
$ws.Range("E17").Value = "  -14.56%  "
$ws.Range("D18").Value = "26.583.59"
$ws.Range("E18").Value = "  -7.25%  "
$ws.Range("D19").Value = "5.006"
$ws.Range("E19").Value = "  -9.18%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "10.78"
$ws.Range("E21").Value = "  -6.11%  "
$ws.Range("D22").Value = "185.52"
$ws.Range("E22").Value = "  -12.16%  "
$ws.Range("D23").Value = "6.201"
$ws.Range("E23").Value = "  -9.80%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "144.66"
$ws.Range("E25").Value = "  -6.29%  "
$ws.Range("D26").Value = "7.507"
$ws.Range("E26").Value = "  -4.84%  "
$ws.Range("D27").Value = "0.1137"
$ws.Range("E27").Value = "  -11.00%  "
$ws.Range("D28").Value = "15.48"
$ws.Range("E28").Value = "  -5.89%  "
$ws.Range("D29").Value = "1.301"
$ws.Range("E29").Value = "  -8.69%  "
$ws.Range("D30").Value = "0.05691"
$ws.Range("E30").Value = "  -8.90%  "
$ws.Range("D31").Value = "1.327"
$ws.Range("E31").Value = "  -6.75%  "
$ws.Range("D32").Value = "3.484"
$ws.Range("E32").Value = "  -8.13%  "
$ws.Range("D33").Value = "3.480"
$ws.Range("E33").Value = "  -7.52%  "
$ws.Range("D34").Value = "1.639"
$ws.Range("E34").Value = "  -5.05%  "
$ws.Range("D35").Value = "1.010"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.5930"
$ws.Range("E36").Value = "  -7.73%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.368"
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("D38").Value = "2.644"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").Value = "0.01592"
$ws.Range("E39").Value = "  -7.22%  "
$ws.Range("D40").Value = "1.071.52"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("D41").Value = "5.907"
$ws.Range("E41").Value = "  -9.07%  "
$ws.Range("D42").Value = "0.8584"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "98.08"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("D45").Value = "1.842.86"
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("D46").Value = "56.41"
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("D47").Value = "0.00000000106"
$ws.Range("E47").Value = "  -5.15%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("D49").Value = "8.020"
$ws.Range("E49").Value = "  -4.15%  "
$ws.Range("D50").Value = "0.4311"
$ws.Range("E50").Value = "  -3.85%  "
$ws.Range("D51").Value = "0.05198"
$ws.Range("E51").Value = "  -4.88%  "

# Reset style on the coerced cells back to Normal (removes the temporary
# text number format while keeping the stored value as text).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
